# Update the three-digit x one-digit multiplication drill numbers.
# Each cell's text is unique in the document, so a targeted
# Find/Replace (MatchWholeWord-ish via exact "NNN×N=" token) on the
# whole-document Range is sufficient and safe — no table structure
# changes are required. Order matters only for the single pair where a
# replacement's new text equals another cell's original text
# (671×6= is produced by row2/cell1 after row2/cell1's own old text
# 671×6= has already been consumed), so we process top-to-bottom,
# left-to-right exactly as the cells appear in the document.

$d = $word.ActiveDocument

# Execute(FindText, MatchCase, MatchWholeWord, MatchWildcards,
#         MatchSoundsLike, MatchAllWordForms, Forward, Wrap, Format,
#         ReplaceWith, Replace)
function Replace-Problem($oldText, $newText) {
    $d.Content.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2) | Out-Null
}

Replace-Problem "867×9=" "116×7="
Replace-Problem "268×7=" "932×4="
Replace-Problem "428×2=" "614×5="
Replace-Problem "606×2=" "581×6="
Replace-Problem "760×2=" "231×6="

Replace-Problem "671×6=" "629×9="
Replace-Problem "361×5=" "887×5="
Replace-Problem "557×2=" "671×6="
Replace-Problem "839×4=" "374×7="
Replace-Problem "857×6=" "261×3="

Replace-Problem "953×5=" "396×7="
Replace-Problem "693×6=" "809×5="
Replace-Problem "528×4=" "429×8="
Replace-Problem "644×4=" "663×9="
Replace-Problem "672×7=" "808×7="

Replace-Problem "266×4=" "295×8="
Replace-Problem "743×8=" "631×6="
Replace-Problem "796×9=" "900×4="
Replace-Problem "338×9=" "101×6="
Replace-Problem "333×5=" "861×8="

Replace-Problem "566×4=" "502×6="
Replace-Problem "549×4=" "160×5="
Replace-Problem "602×2=" "330×5="
Replace-Problem "472×7=" "735×6="
Replace-Problem "497×8=" "866×6="

Write-Host "Updated 25 multiplication problems"
